$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet 1: "Dep or Non-Dep Without LOSOCV"
# ----------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Insert 3 new rows after row 7 (before the second table's title row 10)
# and copy the number-format/border styling from row 7 down into them.
$ws1.Range("A7:F7").Copy()
$ws1.Range("A8:A10").EntireRow.Insert()
$ws1.Range("A8:F10").PasteSpecial(-4122)

# Row 8: Decision Tree results
$ws1.Range("A8").Value = "Decision Tree"
$ws1.Range("B8").Value = 0.89156626506024095
$ws1.Range("C8").Value = 0.85714285714285698
$ws1.Range("D8").Value = 0.88235294117647001
$ws1.Range("E8").Value = 0.86956521739130399
$ws1.Range("F8").Value = 0.89015606242497003

# Row 9: Naive Bayes (no results yet - blank values)
$ws1.Range("A9").Value = "Naïve Bayes"

# Row 10: AdaBoost (no results yet - blank values)
$ws1.Range("A10").Value = "AdaBoost"

# Append 3 new rows at the bottom of the second table (after row 17,
# formerly row 14) for Decision Tree / Naive Bayes / AdaBoost results.
$ws1.Range("A17:F17").Copy()
$ws1.Range("A18:F18").PasteSpecial(-4122)
$ws1.Range("A17:E17").Copy()
$ws1.Range("A19:E20").PasteSpecial(-4122)
$ws1.Range("E17").Copy()
$ws1.Range("F19:F20").PasteSpecial(-4122)

$ws1.Range("A18").Value = "Decision Tree"
$ws1.Range("B18").Value = 0.686746987951807
$ws1.Range("C18").Value = 0.625
$ws1.Range("D18").Value = 0.58823529411764697
$ws1.Range("E18").Value = 0.60606060606060597
$ws1.Range("F18").Value = 0.67166866746698595

$ws1.Range("A19").Value = "Naïve Bayes"
$ws1.Range("A20").Value = "AdaBoost"

$ws1.Application.CutCopyMode = $false

# ----------------------------------------------------------------------
# Sheet 2: "Dep or Non-Dep With LOSOCV"
# ----------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Note: E7 on this sheet carries a special (non-standard) 5-decimal
# number format, so copy A-D/F separately from E to avoid propagating it.
$ws2.Range("A7:D7").Copy()
$ws2.Range("A8:A10").EntireRow.Insert()
$ws2.Range("A8:D10").PasteSpecial(-4122)
$ws2.Range("D7").Copy()
$ws2.Range("E8:E10").PasteSpecial(-4122)
$ws2.Range("F7").Copy()
$ws2.Range("F8:F10").PasteSpecial(-4122)

$ws2.Range("A8").Value = "Decision Tree"
$ws2.Range("B8").Value = 0.79844135662317395
$ws2.Range("C8").Value = 0.381818181818181
$ws2.Range("D8").Value = 0.33851239669421401
$ws2.Range("E8").Value = 0.350375079465988
$ws2.Range("F8").Value = 0.79844135662317395

$ws2.Range("A9").Value = "Naïve Bayes"
$ws2.Range("A10").Value = "AdaBoost"

$ws2.Range("A17:F17").Copy()
$ws2.Range("A18:F20").PasteSpecial(-4122)

$ws2.Range("A18").Value = "Decision Tree"
$ws2.Range("B18").Value = 0.69697211879029997
$ws2.Range("C18").Value = 0.4
$ws2.Range("D18").Value = 0.230578512396694
$ws2.Range("E18").Value = 0.28082200627655102
$ws2.Range("F18").Value = 0.69697211879029997

$ws2.Range("A19").Value = "Naïve Bayes"
$ws2.Range("A20").Value = "AdaBoost"

$ws2.Application.CutCopyMode = $false
